$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "Talbert was a part of a small team ..." paragraph (Work Performance /
# accomplishments section): the text was expanded with a couple of new
# sentences about the team's process, and "in house" was hyphenated to
# "in-house".
# ---------------------------------------------------------------------------
$old1 = "Talbert was a part of a small team that developed and completed an in house reporting application called The Enlightenment Portal. Talbert played a key role in the development and aided the design of the user interface for the application. This application is used by multiple Extended Campuses members; including Blanche Johnson and Kevin Hayes. "
$new1 = "Talbert was a part of a small team that developed and completed an in-house reporting application called The Enlightenment Portal. Talbert played a key role in the development and aided the design of the user interface for the application. This application is used by multiple Extended Campuses members; including Blanche Johnson and Kevin Hayes. The team worked closely with the Business Analyst and each other to plan out the functional and design requirements. In addition, the team implemented a contemporary industry standard software architectural pattern for implementing user interface, called Model-View-Controller (MVC)."

$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null
